$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 269, shifting existing rows 269:353 down to 270:354
$ws.Rows("269:269").Insert()

# Populate the newly inserted row 269 with the new weekly data point
$ws.Range("A269").Value = 7
$ws.Range("B269").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C269").Value = 'Ñuble'
$ws.Range("D269").Value = 44985
$ws.Range("E269").Value = 16
$ws.Range("F269").Value = 100112003
$ws.Range("G269").Value = 'Ajo'
$ws.Range("H269").Value = 'Chino'
$ws.Range("I269").Value = 'Primera'
$ws.Range("J269").Value = 50
$ws.Range("K269").Value = 21000
$ws.Range("L269").Value = 21000
$ws.Range("M269").Value = 21000
$ws.Range("N269").Value = '$/malla 10 kilos'
$ws.Range("O269").Value = 'China'
$ws.Range("P269").Value = 2100
$ws.Range("Q269").Value = 10
$ws.Range("R269").Value = 'Hortaliza'
